$wb = $excel.ActiveWorkbook

# --- Sheet "Human" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Human")

# Update stop-codon related values in column C (rows 12-15 shift up by one codon)
$ws1.Range("C12").Value = "TCA"
$ws1.Range("C13").Value = "TTA"
$ws1.Range("C14").Value = "TTG"
$ws1.Range("C15").ClearContents()

# Column D row 17 gets the value previously held by row 18
$ws1.Range("D17").Value = "TCC"

# Remove now-empty trailing rows 18 and 19 (their only content moved up / was dropped)
$ws1.Range("A18:F19").Delete()

# --- Sheet "Mosquito" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Mosquito")

$ws2.Range("C12").Value = "TCA"
$ws2.Range("C13").Value = "TCT"
$ws2.Range("C14").Value = "TTA"
$ws2.Range("C15").ClearContents()

$ws2.Range("D17").Value = "TCG"

$ws2.Range("A18:F19").Delete()
